$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (this string is shared by Overview!B2, Overview!C2, zh-cn!B2, de-de!B2)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Handback report: for each locale sheet, fill in the columns that get
#    populated once a handback has happened:
#      E2 "Latest Target File"       -> same file/link as the source file (A2)
#      F2 "Latest Handback File"     -> same file/link as the handoff file (C2)
#      G2 "Latest Handback DateTime" -> the handback timestamp
# ---------------------------------------------------------------------------

# Look up the hyperlink already attached to a given cell address on a sheet so
# the new "target"/"handback" links reuse the exact same address & display text.
function Get-HyperlinkByCell($ws, $cellAddress) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellAddress) {
            return $h
        }
    }
    return $null
}

$handbackDateTimes = @{ "zh-cn" = "2016-02-24 08:30:43"; "de-de" = "2016-02-24 08:31:04" }

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $sourceLink = Get-HyperlinkByCell $ws "`$A`$2"
    $handoffLink = Get-HyperlinkByCell $ws "`$C`$2"

    $sourceAddress = $ws.Range("A2").Value
    $sourceDisplay = $ws.Range("A2").Value
    if ($sourceLink -ne $null) {
        $sourceAddress = $sourceLink.Address
        $sourceDisplay = $sourceLink.TextToDisplay
    }

    $handoffAddress = $ws.Range("C2").Value
    $handoffDisplay = $ws.Range("C2").Value
    if ($handoffLink -ne $null) {
        $handoffAddress = $handoffLink.Address
        $handoffDisplay = $handoffLink.TextToDisplay
    }

    $ws.Hyperlinks.Add($ws.Range("E2"), $sourceAddress, "", "", $sourceDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $handoffAddress, "", "", $handoffDisplay) | Out-Null

    $ws.Range("G2").Value = $handbackDateTimes[$sheetName]
}
